# Product Backlog update — "Adding update to backlogs"
#
# Assigns Stas to two previously-unassigned user stories (marking them
# Completed with start/finish dates), assigns Stas to a third story that
# still needs work, and trims Neja off the "integrate database/backend/
# frontend" story so only Kamila remains assigned.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Agile Product Backlog")

# Row 8 — "A user can order the search results based on price": now done by Stas
$ws.Range("E8").Value = "Stas"
$ws.Range("F8").Value = "26th Jan"
$ws.Range("G8").Value = "27th Jan"
$ws.Range("K8").Value = "Completed"

# Row 9 — "A user can order the search results based on distance": now done by Stas
$ws.Range("E9").Value = "Stas"
$ws.Range("F9").Value = "26th Jan"
$ws.Range("G9").Value = "27th Jan"
$ws.Range("K9").Value = "Completed"

# Row 16 — "A user can set a distance range for their search": assigned to Stas, in progress
$ws.Range("E16").Value = "Stas"
$ws.Range("F16").Value = "27th Jan"
$ws.Range("G16").Value = "30th Jan"
$ws.Range("K16").Value = "Completed"

# Row 25 — "As a programmer, I want to integrate the database, backend and
# frontend..." now only assigned to Kamila (Neja dropped off the task)
$ws.Range("E25").Value = "Kamila"
